$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6094.6665
$ws.Range("J62").Value = 4050.5
$ws.Range("L62").Value = 4050.5
$ws.Range("N62").Value = -5298.5
$ws.Range("H65").Value = 6094.6665
$ws.Range("J65").Value = 4050.5
$ws.Range("L65").Value = 20252.5
$ws.Range("N65").Value = -26492.5
$ws.Range("H135").Value = 4855.36
$ws.Range("I135").Value = 1121.5385
$ws.Range("K135").Value = 10093.8465
$ws.Range("M135").Value = -7558.846500000001
$ws.Range("H138").Value = 2799.2554
$ws.Range("I138").Value = 1901.7368
$ws.Range("J138").Value = 3408.2856
$ws.Range("K138").Value = 5705.2104
$ws.Range("L138").Value = 10224.8568
$ws.Range("M138").Value = -565.2103999999999
$ws.Range("N138").Value = -20504.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 191226.78
$ws.Range("I32").Value = 270439.47
$ws.Range("J32").Value = 24000
$ws.Range("K32").Value = 270439.47
$ws.Range("L32").Value = 24000
$ws.Range("M32").Value = -270152.47
$ws.Range("N32").Value = -24574
$ws.Range("H54").Value = 20049
$ws.Range("J54").Value = 20049
$ws.Range("L54").Value = 20049
$ws.Range("N54").Value = -21587
$ws.Range("H61").Value = 1727898
$ws.Range("J61").Value = 12503958
$ws.Range("L61").Value = 12503958
$ws.Range("N61").Value = -12504382
$ws.Range("H74").Value = 980905.2
$ws.Range("I74").Value = 1325471.1
$ws.Range("J74").Value = 16120.667
$ws.Range("K74").Value = 1325471.1
$ws.Range("L74").Value = 16120.667
$ws.Range("M74").Value = -1324597.1
$ws.Range("N74").Value = -17868.667
$ws.Range("H77").Value = 980905.2
$ws.Range("I77").Value = 1325471.1
$ws.Range("J77").Value = 16120.667
$ws.Range("K77").Value = 6627355.5
$ws.Range("L77").Value = 80603.33499999999
$ws.Range("M77").Value = -6622987.5
$ws.Range("N77").Value = -89339.33499999999
$ws.Range("H97").Value = 37038828
$ws.Range("J97").Value = 333336830
$ws.Range("L97").Value = 333336830
$ws.Range("N97").Value = -333337822
$ws.Range("H122").Value = 2282
$ws.Range("I122").Value = 2248.6365
$ws.Range("K122").Value = 6745.9095
$ws.Range("M122").Value = -4295.9095
$ws.Range("H136").Value = 1727898
$ws.Range("J136").Value = 12503958
$ws.Range("L136").Value = 37511874
$ws.Range("N136").Value = -37516974

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 245421.33
$ws.Range("I94").Value = 8435.444
$ws.Range("K94").Value = 8435.444
$ws.Range("M94").Value = -7984.444
$ws.Range("H99").Value = 9657.333000000001
$ws.Range("I99").Value = 13573.5
$ws.Range("K99").Value = 13573.5
$ws.Range("M99").Value = -12075.5
$ws.Range("H134").Value = 3093521.8
$ws.Range("I134").Value = 5033.75
$ws.Range("J134").Value = 16682869
$ws.Range("K134").Value = 15101.25
$ws.Range("L134").Value = 50048607
$ws.Range("M134").Value = -12566.25
$ws.Range("N134").Value = -50053677
$ws.Range("H137").Value = 98999
$ws.Range("J137").Value = 98999
$ws.Range("L137").Value = 98999
$ws.Range("N137").Value = -109199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 36845.57
$ws.Range("J22").Value = 63893.312
$ws.Range("L22").Value = 63893.312
$ws.Range("N22").Value = -64593.312
$ws.Range("H43").Value = 24665
$ws.Range("J43").Value = 24665
$ws.Range("L43").Value = 24665
$ws.Range("N43").Value = -25033
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H101").Value = 24665
$ws.Range("J101").Value = 24665
$ws.Range("L101").Value = 24665
$ws.Range("N101").Value = -31155
$ws.Range("H134").Value = 2436.8928
$ws.Range("I134").Value = 1936.1052
$ws.Range("J134").Value = 3494.111
$ws.Range("K134").Value = 5808.3156
$ws.Range("L134").Value = 10482.333
$ws.Range("M134").Value = -3273.3156
$ws.Range("N134").Value = -15552.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 13789.143
$ws.Range("I96").Value = 11024
$ws.Range("K96").Value = 33072
$ws.Range("M96").Value = -31013
$ws.Range("H131").Value = 4955.5
$ws.Range("J131").Value = 7051.1875
$ws.Range("L131").Value = 21153.5625
$ws.Range("N131").Value = -31233.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7886.232
$ws.Range("I132").Value = 6550.0166
$ws.Range("J132").Value = 16794.334
$ws.Range("K132").Value = 19650.0498
$ws.Range("L132").Value = 50383.00199999999
$ws.Range("M132").Value = -17120.0498
$ws.Range("N132").Value = -55443.00199999999
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 2774.5
$ws.Range("I22").Value = 1800.3334
$ws.Range("J22").Value = 3359
$ws.Range("K22").Value = 1800.3334
$ws.Range("L22").Value = 3359
$ws.Range("M22").Value = -1505.3334
$ws.Range("N22").Value = -3949
$ws.Range("H27").Value = 2774.5
$ws.Range("I27").Value = 1800.3334
$ws.Range("J27").Value = 3359
$ws.Range("K27").Value = 1800.3334
$ws.Range("L27").Value = 3359
$ws.Range("M27").Value = -1693.3334
$ws.Range("N27").Value = -3573
$ws.Range("H82").Value = 1566.8077
$ws.Range("I82").Value = 1410.3478
$ws.Range("K82").Value = 1410.3478
$ws.Range("M82").Value = -1049.3478
$ws.Range("H85").Value = 1566.8077
$ws.Range("I85").Value = 1410.3478
$ws.Range("K85").Value = 1410.3478
$ws.Range("M85").Value = -162.3478
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 1000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 123
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 3000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 2616
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3969804
$ws.Range("I132").Value = 4631015.5
$ws.Range("K132").Value = 13893046.5
$ws.Range("M132").Value = -13890516.5
